$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R data for year 2021
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 18.953297329007047

# Copy style from Q4/Q5 so new cells match formatting of the existing table
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("R5").NumberFormat = "0.0"

# Update the active cell selection
$ws.Range("Q8").Select() | Out-Null
